# Update NATMI ligand-receptor pair output (Vip-Vipr1) with new TPM-derived
# statistics, and drop the now-absent "Resolving-Mac" target-cluster row
# (previously row 6); "Resolving-Mac" becomes the new target cluster for
# what is now the last row (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 6 entirely (its data is superseded / no longer present).
$ws.Rows.Item(6).Delete()

# --- Row 2: Sending cluster ECs -> Target cluster FAPs ---
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 1.030436666666667
$ws.Range("H2").Value = 3.09131
$ws.Range("M2").Value = 0.027509
$ws.Range("N2").Value = 0.082527
$ws.Range("O2").Value = 0.003989820149889837
$ws.Range("P2").Value = 0.003989820149889836
$ws.Range("Q2").Value = 0.02834628226333333
$ws.Range("R2").Value = 0.25511654037
$ws.Range("S2").Value = 0.003989820149889837
$ws.Range("T2").Value = 0.003989820149889836

# --- Row 3: Target cluster -> Inflammatory-Mac ---
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("G3").Value = 1.030436666666667
$ws.Range("H3").Value = 3.09131
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.705026333333334
$ws.Range("N3").Value = 8.115079000000001
$ws.Range("O3").Value = 0.3923286404709717
$ws.Range("P3").Value = 0.3923286404709715
$ws.Range("Q3").Value = 2.787358318165556
$ws.Range("R3").Value = 25.08622486349001
$ws.Range("S3").Value = 0.3923286404709717
$ws.Range("T3").Value = 0.3923286404709715

# --- Row 4: Target cluster -> MuSCs ---
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 1.030436666666667
$ws.Range("H4").Value = 3.09131
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.146307
$ws.Range("N4").Value = 0.438921
$ws.Range("O4").Value = 0.02121991408884119
$ws.Range("P4").Value = 0.02121991408884119
$ws.Range("Q4").Value = 0.15076009739
$ws.Range("R4").Value = 1.35684087651
$ws.Range("S4").Value = 0.02121991408884119
$ws.Range("T4").Value = 0.02121991408884119

# --- Row 5: Target cluster -> Resolving-Mac ---
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 1.030436666666667
$ws.Range("H5").Value = 3.09131
$ws.Range("M5").Value = 4.015954666666667
$ws.Range("N5").Value = 12.047864
$ws.Range("O5").Value = 0.5824616252902973
$ws.Range("P5").Value = 0.5824616252902973
$ws.Range("Q5").Value = 4.138186940204444
$ws.Range("R5").Value = 37.24368246184
$ws.Range("S5").Value = 0.5824616252902973
$ws.Range("T5").Value = 0.5824616252902973
